# Fruta / hortaliza, semanal
# The data rows (2-35) have been reshuffled into a new row order (same 34
# records, different positions). Row 3 keeps its original data. Build the
# mapping of new row -> old row, snapshot all original values first (since
# writes must not clobber values we still need to read), then write the
# snapshotted rows back out in their new order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# newRow -> oldRow
$map = @{
    2  = 11
    3  = 3
    4  = 15
    5  = 18
    6  = 27
    7  = 6
    8  = 17
    9  = 13
    10 = 14
    11 = 26
    12 = 19
    13 = 25
    14 = 12
    15 = 29
    16 = 8
    17 = 23
    18 = 22
    19 = 24
    20 = 10
    21 = 30
    22 = 21
    23 = 32
    24 = 28
    25 = 7
    26 = 31
    27 = 9
    28 = 20
    29 = 2
    30 = 16
    31 = 34
    32 = 5
    33 = 35
    34 = 33
    35 = 4
}

$firstCol = 1   # A
$lastCol  = 18  # R

# Snapshot every data row (2..35) across columns A..R before writing anything.
$snapshot = @{}
for ($r = 2; $r -le 35; $r++) {
    $rowVals = @{}
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# Write the snapshotted data back out into the new row order.
foreach ($newRow in $map.Keys) {
    $oldRow = $map[$newRow]
    $rowVals = $snapshot[$oldRow]
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $ws.Cells.Item($newRow, $c).Value2 = $rowVals[$c]
    }
}
